# PROS-13942-JnJANZ - Exclude Sub Category-Update Config Template
#
# The KPI Exclusions Template gains a new "sub_category" exclusion column
# pair: header "sub_category" in D2 (the exclude-field name) and value
# "Out of Scope" in E2 (the value to exclude), mirroring the existing
# A2/B2/C2 exclusion rule on the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "sub_category"
$ws.Range("E2").Value = "Out of Scope"

# Leave the user's cursor where they finished working.
$ws.Range("F24").Select()
